$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the "Female order" column group (row 11)
$ws.Range("E11").Value = "Female order"

# Row 12: Soil Temperature/Moisture Sensor - SHT10 / Red / 5V, plus female-order columns
$ws.Range("B12").Value = "Soil Temperature/Moisture Sensor - SHT10"
$ws.Range("C12").Value = "Red"
$ws.Range("D12").Value = "5V"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "red"

# Row 13: Green / GND
$ws.Range("C13").Value = "Green"
$ws.Range("D13").Value = "GND"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = "black"

# Row 14: Yellow (clockPin) / D11
$ws.Range("C14").Value = "Yellow  (clockPin)"
$ws.Range("D14").Value = "D11"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = "Yellow"

# Row 15: Blue (dataPin) / D10
$ws.Range("C15").Value = "Blue  (dataPin)"
$ws.Range("D15").Value = "D10"
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "White"

# Widen column B slightly to fit new content (best-fit after the longer entries)
$ws.Columns("B:B").ColumnWidth = 36

# Update the view's scroll position / selection
$excel.ActiveWindow.TopLeftCell = $ws.Range("A8")
$ws.Range("F16").Select() | Out-Null
